$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "37.747.94"
$ws.Range("E2").Value = "  -0.23%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.039.87"
$ws.Range("E3").Value = "  +0.28%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "227.38"
$ws.Range("E5").Value = "  -0.19%  "
$ws.Range("E6").Value = "  -1.13%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "59.56"
$ws.Range("E7").Value = "  -1.33%  "
$ws.Range("E8").Value = "  -0.04%  "
$ws.Range("E9").Value = "  -2.71%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0841"
$ws.Range("E10").Value = "  +2.89%  "
$ws.Range("E11").Value = "  -0.27%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "2.341.55"
$ws.Range("E12").Value = "  +0.33%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "14.46"
$ws.Range("E13").Value = "  -1.61%  "
$ws.Range("E14").Value = "  -0.70%  "
$ws.Range("E15").Value = "  +4.07%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.773"
$ws.Range("E16").Value = "  +2.09%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.040.79"
$ws.Range("E17").Value = "  +0.30%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "37.667.83"
$ws.Range("E18").Value = "  -0.38%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "69.48"
$ws.Range("E19").Value = "  -0.57%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "5.90"
$ws.Range("E20").Value = "  -2.19%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0₃0823"
$ws.Range("E21").Value = "  -0.13%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "223.15"
$ws.Range("E22").Value = "  -1.23%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.43"
$ws.Range("E24").Value = "  +1.20%  "
$ws.Range("E25").Value = "  +2.09%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "168.30"
$ws.Range("E26").Value = "  +2.05%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.37"
$ws.Range("E27").Value = "  +1.00%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.129"
$ws.Range("E28").Value = "  -0.93%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "18.80"
$ws.Range("E29").Value = "  -0.60%  "
$ws.Range("E30").Value = "  -0.39%  "
$ws.Range("E31").Value = "  -0.88%  "
$ws.Range("E32").Value = "  +8.40%  "
$ws.Range("E33").Value = "  -1.43%  "
$ws.Range("E34").Value = "  +0.38%  "
$ws.Range("E35").Value = "  +0.24%  "
$ws.Range("E36").Value = "  +1.27%  "
$ws.Range("E37").Value = "  +3.10%  "
$ws.Range("E38").Value = "  +4.47%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.00"
$ws.Range("E39").Value = "  -0.04%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "18.06"
$ws.Range("E40").Value = "  +7.77%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.528.41"
$ws.Range("E41").Value = "  -0.59%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "97.27"
$ws.Range("E42").Value = "  +0.27%  "
$ws.Range("E43").Value = "  -1.18%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.84"
$ws.Range("E44").Value = "  +1.05%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "4.23"
$ws.Range("E45").Value = "  +6.74%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0905"
$ws.Range("E46").Value = "  -1.81%  "
$ws.Range("E47").Value = "  -0.33%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.01"
$ws.Range("E48").Value = "  -0.07%  "
$ws.Range("E49").Value = "  -0.73%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "6.99"
$ws.Range("E50").Value = "  -1.93%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.230.18"
$ws.Range("E51").Value = "  +0.31%  "
